$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 480.36365
$ws.Range("I5").Value = 480.36365
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 480.36365
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -365.36365
$ws.Range("H70").Value = 1801770.5
$ws.Range("I70").Value = 8990929
$ws.Range("J70").Value = 4480.875
$ws.Range("K70").Value = 26972787
$ws.Range("L70").Value = 13442.625
$ws.Range("M70").Value = -26972517
$ws.Range("N70").Value = -13982.625
$ws.Range("H73").Value = 1801770.5
$ws.Range("I73").Value = 8990929
$ws.Range("J73").Value = 4480.875
$ws.Range("K73").Value = 26972787
$ws.Range("L73").Value = 13442.625
$ws.Range("M73").Value = -26971851
$ws.Range("N73").Value = -15314.625
$ws.Range("H80").Value = 1263829.6
$ws.Range("I80").Value = 2526030
$ws.Range("J80").Value = 1629.4445
$ws.Range("K80").Value = 7578090
$ws.Range("L80").Value = 4888.333500000001
$ws.Range("M80").Value = -7577092
$ws.Range("N80").Value = -6884.333500000001
$ws.Range("H83").Value = 1263829.6
$ws.Range("I83").Value = 2526030
$ws.Range("J83").Value = 1629.4445
$ws.Range("K83").Value = 22734270
$ws.Range("L83").Value = 14665.0005
$ws.Range("M83").Value = -22729278
$ws.Range("N83").Value = -24649.0005
$ws.Range("H86").Value = 15428537
$ws.Range("I86").Value = 9000
$ws.Range("J86").Value = 20054398
$ws.Range("K86").Value = 9000
$ws.Range("L86").Value = 20054398
$ws.Range("M86").Value = -7877
$ws.Range("N86").Value = -20056644
$ws.Range("H89").Value = 15428537
$ws.Range("I89").Value = 9000
$ws.Range("J89").Value = 20054398
$ws.Range("K89").Value = 45000
$ws.Range("L89").Value = 100271990
$ws.Range("M89").Value = -39384
$ws.Range("N89").Value = -100283222
$ws.Range("H97").Value = 2999.8
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2999.8
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 8999.400000000001
$ws.Range("N97").Value = -9991.400000000001
$ws.Range("H98").Value = 3947.8
$ws.Range("I98").Value = 4219.778
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 4219.778
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -2721.778
$ws.Range("N98").Value = -4496
$ws.Range("H112").Value = 54796.633
$ws.Range("I112").Value = 252249.75
$ws.Range("J112").Value = 2142.4666
$ws.Range("K112").Value = 756749.25
$ws.Range("L112").Value = 6427.399800000001
$ws.Range("M112").Value = -755641.25
$ws.Range("N112").Value = -8643.399800000001
$ws.Range("H122").Value = 3947.8
$ws.Range("I122").Value = 4219.778
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 12659.334
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -10209.334
$ws.Range("N122").Value = -9400

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1888.8334
$ws.Range("I61").Value = 1333.5
$ws.Range("J61").Value = 2999.5
$ws.Range("K61").Value = 1333.5
$ws.Range("L61").Value = 2999.5
$ws.Range("M61").Value = -1121.5
$ws.Range("N61").Value = -3423.5
$ws.Range("H88").Value = 10418945
$ws.Range("I88").Value = 23810982
$ws.Range("J88").Value = 2916.2222
$ws.Range("K88").Value = 23810982
$ws.Range("L88").Value = 2916.2222
$ws.Range("M88").Value = -23810576
$ws.Range("N88").Value = -3728.2222
$ws.Range("H91").Value = 10418945
$ws.Range("I91").Value = 23810982
$ws.Range("J91").Value = 2916.2222
$ws.Range("K91").Value = 23810982
$ws.Range("L91").Value = 2916.2222
$ws.Range("M91").Value = -23809578
$ws.Range("N91").Value = -5724.2222
$ws.Range("H102").Value = 8396.444
$ws.Range("I102").Value = 2224.1428
$ws.Range("J102").Value = 29999.5
$ws.Range("K102").Value = 2224.1428
$ws.Range("L102").Value = 29999.5
$ws.Range("M102").Value = -602.1428000000001
$ws.Range("N102").Value = -33243.5
$ws.Range("H132").Value = 3874.4546
$ws.Range("I132").Value = 3813.1428
$ws.Range("J132").Value = 3981.75
$ws.Range("K132").Value = 11439.4284
$ws.Range("L132").Value = 11945.25
$ws.Range("M132").Value = -8909.428400000001
$ws.Range("N132").Value = -17005.25
$ws.Range("H136").Value = 1888.8334
$ws.Range("I136").Value = 1333.5
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 4000.5
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -1450.5
$ws.Range("N136").Value = -14098.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 33335486
$ws.Range("I86").Value = 50001628
$ws.Range("J86").Value = 3201.8
$ws.Range("K86").Value = 50001628
$ws.Range("L86").Value = 3201.8
$ws.Range("M86").Value = -50000505
$ws.Range("N86").Value = -5447.8
$ws.Range("H89").Value = 33335486
$ws.Range("I89").Value = 50001628
$ws.Range("J89").Value = 3201.8
$ws.Range("K89").Value = 250008140
$ws.Range("L89").Value = 16009
$ws.Range("M89").Value = -250002524
$ws.Range("N89").Value = -27241
$ws.Range("H105").Value = 2392.5908
$ws.Range("I105").Value = 2311
$ws.Range("J105").Value = 2610.1667
$ws.Range("K105").Value = 2311
$ws.Range("L105").Value = 2610.1667
$ws.Range("M105").Value = -564
$ws.Range("N105").Value = -6104.1667

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2220.524
$ws.Range("I58").Value = 1921
$ws.Range("J58").Value = 2291
$ws.Range("K58").Value = 1921
$ws.Range("L58").Value = 2291
$ws.Range("M58").Value = -1718
$ws.Range("N58").Value = -2697
$ws.Range("H132").Value = 3597.8518
$ws.Range("I132").Value = 2970.0625
$ws.Range("J132").Value = 4511
$ws.Range("K132").Value = 8910.1875
$ws.Range("L132").Value = 13533
$ws.Range("M132").Value = -6380.1875
$ws.Range("N132").Value = -18593
$ws.Range("H134").Value = 2499.95
$ws.Range("I134").Value = 2024.6818
$ws.Range("J134").Value = 3080.8333
$ws.Range("K134").Value = 6074.0454
$ws.Range("L134").Value = 9242.499899999999
$ws.Range("M134").Value = -3539.0454
$ws.Range("N134").Value = -14312.4999
$ws.Range("H136").Value = 2220.524
$ws.Range("I136").Value = 1921
$ws.Range("J136").Value = 2291
$ws.Range("K136").Value = 5763
$ws.Range("L136").Value = 6873
$ws.Range("M136").Value = -3213
$ws.Range("N136").Value = -11973

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 268476.06
$ws.Range("I68").Value = 1001212.25
$ws.Range("J68").Value = 2026.5454
$ws.Range("K68").Value = 3003636.75
$ws.Range("L68").Value = 6079.6362
$ws.Range("M68").Value = -3002825.75
$ws.Range("N68").Value = -7701.6362
$ws.Range("H71").Value = 268476.06
$ws.Range("I71").Value = 1001212.25
$ws.Range("J71").Value = 2026.5454
$ws.Range("K71").Value = 9010910.25
$ws.Range("L71").Value = 18238.9086
$ws.Range("M71").Value = -9006854.25
$ws.Range("N71").Value = -26350.9086

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2411.2856
$ws.Range("I43").Value = 2411.2856
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 2411.2856
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -2260.2856
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H70").Value = 11864.75
$ws.Range("I70").Value = 11418.429
$ws.Range("J70").Value = 14989
$ws.Range("K70").Value = 11418.429
$ws.Range("L70").Value = 14989
$ws.Range("M70").Value = -11148.429
$ws.Range("N70").Value = -15529
$ws.Range("H73").Value = 11864.75
$ws.Range("I73").Value = 11418.429
$ws.Range("J73").Value = 14989
$ws.Range("K73").Value = 11418.429
$ws.Range("L73").Value = 14989
$ws.Range("M73").Value = -10482.429
$ws.Range("N73").Value = -16861
$ws.Range("H122").Value = 3247.4443
$ws.Range("I122").Value = 2769.4546
$ws.Range("J122").Value = 3998.5715
$ws.Range("K122").Value = 8308.363799999999
$ws.Range("L122").Value = 11995.7145
$ws.Range("M122").Value = -5858.363799999999
$ws.Range("N122").Value = -16895.7145
$ws.Range("H126").Value = 6623.353
$ws.Range("I126").Value = 9027.111000000001
$ws.Range("J126").Value = 3919.125
$ws.Range("K126").Value = 27081.333
$ws.Range("L126").Value = 11757.375
$ws.Range("M126").Value = -24611.333
$ws.Range("N126").Value = -16697.375
$ws.Range("H132").Value = 439782.4
$ws.Range("I132").Value = 670867.5600000001
$ws.Range("J132").Value = 6497.75
$ws.Range("K132").Value = 2012602.68
$ws.Range("L132").Value = 19493.25
$ws.Range("M132").Value = -2010072.68
$ws.Range("N132").Value = -24553.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2357.0356
$ws.Range("I16").Value = 2080.625
$ws.Range("J16").Value = 2725.5833
$ws.Range("K16").Value = 2080.625
$ws.Range("L16").Value = 2725.5833
$ws.Range("M16").Value = -1910.625
$ws.Range("N16").Value = -3065.5833
$ws.Range("H55").Value = 648.5
$ws.Range("I55").Value = 775.8333
$ws.Range("J55").Value = 457.5
$ws.Range("K55").Value = 775.8333
$ws.Range("L55").Value = 457.5
$ws.Range("M55").Value = -602.8333
$ws.Range("N55").Value = -803.5
$ws.Range("H93").Value = 1446.2858
$ws.Range("I93").Value = 1474.037
$ws.Range("J93").Value = 1352.625
$ws.Range("K93").Value = 1474.037
$ws.Range("L93").Value = 1352.625
$ws.Range("M93").Value = -226.037
$ws.Range("N93").Value = -3848.625
$ws.Range("H100").Value = 999
$ws.Range("I100").Value = 999
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 999
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -458
$ws.Range("H122").Value = 3416.9268
$ws.Range("I122").Value = 2659.95
$ws.Range("J122").Value = 4137.857
$ws.Range("K122").Value = 7979.849999999999
$ws.Range("L122").Value = 12413.571
$ws.Range("M122").Value = -5529.849999999999
$ws.Range("N122").Value = -17313.571
$ws.Range("H132").Value = 4868.0933
$ws.Range("I132").Value = 3203.5557
$ws.Range("J132").Value = 7677
$ws.Range("K132").Value = 9610.667099999999
$ws.Range("L132").Value = 23031
$ws.Range("M132").Value = -7080.667099999999
$ws.Range("N132").Value = -28091
$ws.Range("H136").Value = 7083.5
$ws.Range("I136").Value = 3070.5557
$ws.Range("J136").Value = 14306.8
$ws.Range("K136").Value = 9211.667099999999
$ws.Range("L136").Value = 42920.39999999999
$ws.Range("M136").Value = -6661.667099999999
$ws.Range("N136").Value = -48020.39999999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10531224
$ws.Range("I81").Value = 2032.4445
$ws.Range("J81").Value = 20007496
$ws.Range("K81").Value = 4064.889
$ws.Range("L81").Value = 40014992
$ws.Range("M81").Value = -3003.889
$ws.Range("N81").Value = -40017114
$ws.Range("H84").Value = 10531224
$ws.Range("I84").Value = 2032.4445
$ws.Range("J84").Value = 20007496
$ws.Range("K84").Value = 20324.445
$ws.Range("L84").Value = 200074960
$ws.Range("M84").Value = -15020.445
$ws.Range("N84").Value = -200085568
$ws.Range("H126").Value = 2036.3077
$ws.Range("I126").Value = 1968.5
$ws.Range("J126").Value = 2262.3333
$ws.Range("K126").Value = 5905.5
$ws.Range("L126").Value = 6786.999899999999
$ws.Range("M126").Value = -3435.5
$ws.Range("N126").Value = -11726.9999
$ws.Range("H132").Value = 276551.03
$ws.Range("I132").Value = 346956.25
$ws.Range("J132").Value = 4317.6
$ws.Range("K132").Value = 1040868.75
$ws.Range("L132").Value = 12952.8
$ws.Range("M132").Value = -1038338.75
$ws.Range("N132").Value = -18012.8
